# Auto-generated from the author diff: updates crypto price/volume table.
# D-column price values that could be misread as numbers are entered with a
# leading apostrophe (forces text, exactly like typing in Excel), then the
# cell style is reset to Normal so no stray quote-prefix style sticks around.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.753.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'2.675.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'600.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'155.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'5.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'0.396"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'29.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "'0.0000194"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "'3.159.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "'65.589.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'2.679.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'12.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "'350.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'70.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'9.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000110"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").Value = "'1.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'0.168"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").Value = "'8.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'535.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("D34").Value = "'6.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").Value = "'5.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").Value = "'0.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").Value = "'20.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "'160.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'42.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'165.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "'4.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "'0.0615"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'22.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.648"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0261"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'20.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
